$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(367, 44441, 2, 12, 94.22850412249706),
    @(368, 44442, 1, 9, 70.67137809187278),
    @(369, 44443, 3, 9, 70.67137809187278),
    @(370, 44444, 1, 9, 70.67137809187278),
    @(371, 44445, 6, 13, 102.0808794660385),
    @(372, 44446, 5, 18, 141.3427561837456),
    @(373, 44447, 0, 18, 141.3427561837456),
    @(374, 44448, 0, 16, 125.6380054966627)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

# Column A on the new rows uses the same date style as the rest of the
# column (centered, bold, bordered, custom datetime number format), so
# copy that formatting down from the last existing row (366).
$ws.Cells.Item(366, 1).Copy()
$ws.Range("A367:A374").PasteSpecial(-4122)
$excel.CutCopyMode = $false
